# The source workbook's "GerdauHP" sheet had a duplicate/erroneous profile
# row (old row 24, profile "310x107") that needs to be removed. Deleting it
# shifts all subsequent rows up by one (old row 25 becomes new row 24, etc.),
# shrinking the used range from A1:L31 down to A1:L30 and dropping the
# trailing row that previously sat at r="31".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GerdauHP")
$ws.Rows.Item(24).Delete()

# Leave the cursor where the author left it after the edit.
$ws.Range("D32").Select()
